$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values
$ws.Range("B2").Value = 5.5
$ws.Range("C5").Value = 14

# Update selection to B2
$ws.Range("B2").Select()
